$wb = $excel.ActiveWorkbook

# --- Login sheet: password reset from "Admin1" to "Admin" ---
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("A1").Value = "Admin"

# --- Emp sheet: new columns (E..J) and two additional rows of data ---
$wsEmp = $wb.Worksheets.Item("Emp")

# Shift the existing 4 columns (Admin/admin123/Aswini101/Selenium/101) two
# columns to the right (C:F) and put two brand-new, unbordered columns
# (Admin/admin123) in front at A:B.
$wsEmp.Range("F1").Value = 101
$wsEmp.Range("E1").Value = "Ram"
$wsEmp.Range("D1").Value = "Selenium"
$wsEmp.Range("C1").Value = "Aswini101"
$wsEmp.Range("B1").Value = "admin123"
$wsEmp.Range("A1").Value = "Admin"

# E1:F1 need the same thin-all-around border already used by C1:D1 - copy
# that formatting across instead of re-creating it from scratch.
$wsEmp.Range("C1:D1").Copy()
$wsEmp.Range("E1:F1").PasteSpecial(-4122)

# A1:B1 must end up with NO border (plain default style).
$wsEmp.Range("A1:B1").Borders.LineStyle = -4142

# G1:J1 - new cells with a left+right thin border only.
$wsEmp.Range("G1").Value = "Capture.png"
$wsEmp.Range("H1").Value = "user1"
$wsEmp.Range("I1").Value = "password"
$wsEmp.Range("J1").Value = "enabled"
foreach ($col in @("G", "H", "I", "J")) {
    $cell = $wsEmp.Range($col + "1")
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# Row 2
$wsEmp.Range("A2").Value = "Admin"
$wsEmp.Range("B2").Value = "admin123"
$wsEmp.Range("C2").Value = "Aswini102"
$wsEmp.Range("D2").Value = "Selenium"
$wsEmp.Range("E2").Value = "Ram"
$wsEmp.Range("F2").Value = 101
$wsEmp.Range("G2").Value = "Capture.png"
$wsEmp.Range("H2").Value = "user1"
$wsEmp.Range("I2").Value = "password"
$wsEmp.Range("J2").Value = "enabled"

$wsEmp.Range("C1:F1").Copy()
$wsEmp.Range("C2:F2").PasteSpecial(-4122)
$wsEmp.Range("G1:J1").Copy()
$wsEmp.Range("G2:J2").PasteSpecial(-4122)

# Row 3
$wsEmp.Range("A3").Value = "Admin"
$wsEmp.Range("B3").Value = "admin123"
$wsEmp.Range("C3").Value = "Aswini103"
$wsEmp.Range("D3").Value = "Selenium"
$wsEmp.Range("E3").Value = "Ram"
$wsEmp.Range("F3").Value = 102
$wsEmp.Range("G3").Value = "Capture.png"
$wsEmp.Range("H3").Value = "user2"
$wsEmp.Range("I3").Value = "password"
$wsEmp.Range("J3").Value = "enabled"

$wsEmp.Range("C1:F1").Copy()
$wsEmp.Range("C3:F3").PasteSpecial(-4122)
$wsEmp.Range("G1:J1").Copy()
$wsEmp.Range("G3:J3").PasteSpecial(-4122)

# --- Selections / active sheet ---
$wsLogin.Range("A4").Select()
$wsEmp.Activate()
$wsEmp.Range("B4").Select()
